# Implementeer geavanceerd filtersysteem
# - Voegt een nieuwe kolom "OpdrachtType" toe (kolom C), bestaande kolommen C:F schuiven op naar D:G.
# - Vult OpdrachtType in voor de bestaande 12 opdrachten.
# - Voegt nieuwe rijen (13-21) toe met "test" categorie opdrachten en types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Nieuwe kolom C (OpdrachtType) invoegen; bestaande C..F schuiven naar D..G ---
$ws.Columns.Item(3).Insert()

# --- 2. Kopteksten (rij 1) ---
$ws.Cells.Item(1,3).Value = "OpdrachtType"

# --- 3. OpdrachtType vullen voor bestaande opdrachten (rijen 2-12) ---
$opdrachtTypes = @{
  2  = "feitenkennis"
  3  = "toepassing"
  4  = "toepassing"
  5  = "toepassing"
  6  = "toepassing"
  7  = "feitenkennis"
  8  = "Tekenen"
  9  = "Tekenen"
  10 = "Tekenen"
  11 = "Tekenen"
  12 = "Communicatie"
}
foreach ($r in $opdrachtTypes.Keys) {
  $ws.Cells.Item($r,3).Value = $opdrachtTypes[$r]
}

# --- 4. Nieuwe rijen 13-21 toevoegen (test-categorie en CriteriaModal voorbeelddata) ---
$newRows = @(
  @{ Row=13; A="Algemeen"; B="test"; C="Communicatie"; D="test vraag" }
  @{ Row=14; B="test"; C="Feitenkennis"; D="Wat zijn de drie fases van weefselherstel?" }
  @{ Row=15; B="test"; C="Begrijpen"; D="Waarom is de ontstekingsfase belangrijk voor genezing?" }
  @{ Row=16; B="test"; C="Toepassing"; D="Demonstreer een effectieve oefening voor het versterken van de quadriceps." }
  @{ Row=17; B="test"; C="Uitleggen"; D="Leg aan een patiënt uit wat het verschil is tussen artritis en artrose." }
  @{ Row=18; B="wel"; D="Teken de anatomie van het kniegewricht, inclusief de kruisbanden en menisci." }
  @{ Row=19; B="test"; C="Communicatie"; D="Voer een rollenspel uit waarin je een patiënt motiveert om zijn thuisoefeningen te doen." }
  @{ Row=20; B="test"; C="Fysiotherapie"; D="Welke manuele techniek zou je toepassen bij een patiënt met een tenniselleboog?" }
  @{ Row=21; B="test"; C="Praktijk"; D="Analyseer deze (fictieve) casus en stel een behandelplan op voor de eerste twee weken." }
)

foreach ($rowData in $newRows) {
  $r = $rowData.Row
  if ($rowData.ContainsKey("A")) { $ws.Cells.Item($r,1).Value = $rowData.A }
  if ($rowData.ContainsKey("B")) { $ws.Cells.Item($r,2).Value = $rowData.B }
  if ($rowData.ContainsKey("C")) { $ws.Cells.Item($r,3).Value = $rowData.C }
  if ($rowData.ContainsKey("D")) { $ws.Cells.Item($r,4).Value = $rowData.D }
}

# --- 5. Kolombreedtes instellen (best effort benadering van de opgeslagen bestFit breedtes) ---
$ws.Columns.Item(1).ColumnWidth = 12.963541666666666
$ws.Columns.Item(2).ColumnWidth = 12.264322916666666
$ws.Columns.Item(3).ColumnWidth = 12.264322916666666
$ws.Columns.Item(4).ColumnWidth = 67.76432291666667
$ws.Columns.Item(5).ColumnWidth = 53.264322916666664
$ws.Columns.Item(6).ColumnWidth = 13.166666666666666
$ws.Columns.Item(7).ColumnWidth = 18.264322916666668

# --- 6. Selectie en scrollpositie bijwerken ---
[void]$ws.Range("B18").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
